$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.455.48'
$ws.Range('E2').Value = '  +5.56%  '

$ws.Range('D3').Value = '3.617.08'
$ws.Range('E3').Value = '  +5.09%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.74%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '191.45'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.90%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.645'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.06%  '

$ws.Range('D8').Value = '3.608.49'
$ws.Range('E8').Value = '  +5.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.07%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.179'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.50%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.666'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.49%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.84'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.85%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000290'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.30%  '

$ws.Range('D15').Value = '4.186.56'
$ws.Range('E15').Value = '  +4.94%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.96%  '

$ws.Range('D17').Value = '3.609.60'
$ws.Range('E17').Value = '  +4.75%  '

$ws.Range('D18').Value = '70.294.67'
$ws.Range('E18').Value = '  +5.48%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.63%  '

$ws.Range('E20').Value = '  +0.71%  '

$ws.Range('E21').Value = '  +4.46%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '493.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +16.79%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.17%  '

$ws.Range('E25').Value = '  +3.75%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '91.09'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.42%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.20%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.76%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.55%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.44%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.73'
$ws.Range('D31').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '638.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.42%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.35'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.53%  '

$ws.Range('E34').Value = '  +7.08%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '65.87'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.11%  '

$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0825'
$ws.Range('E36').Value = '  +7.00%  '

$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.85%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.406'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.14%  '

$ws.Range('E39').Value = '  -0.02%  '

$ws.Range('E40').Value = '  -0.73%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.56'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.04%  '

$ws.Range('D42').Value = '3.307.63'
$ws.Range('E42').Value = '  +3.83%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.37%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.75'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.33%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0453'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.64%  '

$ws.Range('E46').Value = '  +2.73%  '

$ws.Range('E47').Value = '  +2.31%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.13'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.74%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.11%  '

$ws.Range('E50').Value = '  +5.05%  '

$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.04%  '
